# ConfigurablePricing_TestData.xlsx — "Added Filter for Customer in the Discount Feature"
#
# Renames "Pricing" -> "Configurable Pricing" in several labels across the
# Pricing / Product / Customer / Order sheets, and bumps the Web Data
# identifier used on every sheet.

$wb = $excel.ActiveWorkbook

$pricing  = $wb.Worksheets.Item("Pricing")
$product  = $wb.Worksheets.Item("Product")
$customer = $wb.Worksheets.Item("Customer")
$order    = $wb.Worksheets.Item("Order")

# ---------------------------------------------------------------------
# 1. "Web Data 3" -> "Web Data 5" (shared by C1 on every sheet)
# ---------------------------------------------------------------------
$pricing.Range("C1").Value  = "Web Data 5"
$product.Range("C1").Value  = "Web Data 5"
$customer.Range("C1").Value = "Web Data 5"
$order.Range("C1").Value    = "Web Data 5"

# ---------------------------------------------------------------------
# 2. Product sheet labels
# ---------------------------------------------------------------------
$product.Range("D1").Value = "Configurable Pricing category"

$planCat = $product.Range("D2")
$planCat.Value = "Plan Configurable Pricing Category"
$planCat.Characters(6, 12).Font.Name  = "Calibri"
$planCat.Characters(6, 12).Font.Size  = 12
$planCat.Characters(6, 12).Font.Color = 0

# ---------------------------------------------------------------------
# 3. Customer sheet labels (Tiered / Volume Pricing Customer)
# ---------------------------------------------------------------------
$tieredCust = $customer.Range("E1")
$tieredCust.Value = "Tiered Configurable Pricing Customer"
$tieredCust.Characters(8, 13).Font.Name  = "Calibri"
$tieredCust.Characters(8, 13).Font.Size  = 12
$tieredCust.Characters(8, 13).Font.Color = 0

$volumeCust = $customer.Range("E2")
$volumeCust.Value = "Volume Configurable Pricing Customer"
$volumeCust.Characters(8, 13).Font.Name  = "Calibri"
$volumeCust.Characters(8, 13).Font.Size  = 12
$volumeCust.Characters(8, 13).Font.Color = 0

$customer.Rows.Item(1).RowHeight = 30
$customer.Rows.Item(2).RowHeight = 30

# ---------------------------------------------------------------------
# 4. Order sheet mirrors the same two Customer labels
# ---------------------------------------------------------------------
$tieredOrder = $order.Range("D1")
$tieredOrder.Value = "Tiered Configurable Pricing Customer"
$tieredOrder.Characters(8, 13).Font.Name  = "Calibri"
$tieredOrder.Characters(8, 13).Font.Size  = 12
$tieredOrder.Characters(8, 13).Font.Color = 0

$volumeOrder = $order.Range("D2")
$volumeOrder.Value = "Volume Configurable Pricing Customer"
$volumeOrder.Characters(8, 13).Font.Name  = "Calibri"
$volumeOrder.Characters(8, 13).Font.Size  = 12
$volumeOrder.Characters(8, 13).Font.Color = 0

$order.Rows.Item(1).RowHeight = 27.25
$order.Rows.Item(2).RowHeight = 29.3

# ---------------------------------------------------------------------
# 5. Product sheet row 1 height nudges slightly after the label edit
# ---------------------------------------------------------------------
$product.Rows.Item(1).RowHeight = 29.85

# ---------------------------------------------------------------------
# 6. Selection bookkeeping — restore focus on each sheet, Pricing last
#    so it stays the active tab (matching the original workbook state).
# ---------------------------------------------------------------------
[void]$product.Range("C1").Select()
[void]$customer.Range("C1").Select()
[void]$order.Range("C1").Select()
[void]$pricing.Range("C2").Select()
